# Error Calculations and Plots
# Apply missing-data edits to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level corrections on rows that keep their row number ---
# RM 125 (row 19): column D (header "C") now has a value
$ws.Range("D19").Value = -15.5

# RM 135 (row 21): column D (header "C") becomes missing
$ws.Range("D21").ClearContents()

# RM 140 (row 23): column D (header "C") now has a value
$ws.Range("D23").Value = -13.9

# --- Remove two whole rows (RM 232 and SC 92) ---
# Delete the lower row first so the row number of the upper one
# (row 26, "RM 232") is unaffected by the shift.
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# --- Cell-level corrections on the rows that shifted up after the deletions ---
# SC 5 (now row 26): column C (header "B") value becomes missing
$ws.Range("C26").ClearContents()

# SC 101 (now row 27): column C (header "B") gains a value, column D (header "C") becomes missing
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()

# SC 119 (now row 29): column C (header "B") value becomes missing
$ws.Range("C29").ClearContents()

# SC 232 (now row 33): column D (header "C") gains a value
$ws.Range("D33").Value = -14.1
